# Update Data by bot, scripted by HH
# Applies the 688698 balance-sheet refresh to row 2 of the active sheet:
# DATE_TYPE_CODE flips from the "002" (half-year) code to "001" (annual),
# REPORT_DATE moves from 2020-06-30 to 2019-12-31, and every financial
# figure / ratio in O2:AG2 is refreshed to match the new reporting period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J2 (DATE_TYPE_CODE): text "002" -> "001" ----------------------------
# A bare Value assignment of a zero-padded numeric string gets coerced to a
# number by Excel's normal "General" auto-typing, so round-trip the literal
# through a text formula and flatten it back to a plain value with
# PasteSpecial(xlPasteValues). This preserves the inline/shared-string cell
# type and the existing (default) cell style, matching the source data.
$ws.Range("J2").Formula = "=""001"""
$ws.Range("J2").Copy($ws.Range("J2"))
$ws.Range("J2").PasteSpecial(-4163)

# --- N2 (REPORT_DATE): plain text, no numeric coercion risk --------------
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# --- O2:AG2 numeric figures & ratios --------------------------------------
$ws.Range("O2").Value = 451417565.64
$ws.Range("P2").Value = 88828807.73
$ws.Range("Q2").Value = 62765070.75
$ws.Range("R2").Value = 26.7371598747
$ws.Range("S2").Value = 107011127.02
$ws.Range("T2").Value = 2.2338454294
$ws.Range("U2").Value = 56306947.45
$ws.Range("V2").Value = -24.6854814242
$ws.Range("W2").Value = 185976922.11
$ws.Range("X2").Value = 59234015.95
$ws.Range("Y2").Value = -5.3767821701
$ws.Range("Z2").Value = 4472557.95
$ws.Range("AA2").Value = -5.6361506744
$ws.Range("AB2").Value = 265440643.53
$ws.Range("AC2").Value = 30.1298216791
$ws.Range("AD2").Value = 8.3508638931
$ws.Range("AE2").Value = -12.5408424924
$ws.Range("AF2").Value = 216.3490061109
$ws.Range("AG2").Value = 41.1984238687
